$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.01222129826550129, 0.9966874999999999, 0.03400449777442152, 0.9905833333333334, 488.926047, 86.752)
    3  = @(0.01330745655684708, 0.9962916666666666, 0.03533847363175308, 0.9903333333333333, 486.714349, 87.578)
    4  = @(0.01349757201275821, 0.9960624999999999, 0.03955140594994933, 0.9894722222222222, 481.579698, 87.358)
    5  = @(0.01342514414958926, 0.9959583333333333, 0.04085397449813036, 0.9892083333333334, 480.298261, 87.39700000000001)
    6  = @(0.01298716859576598, 0.9961208333333333, 0.04054353564209103, 0.9893666666666666, 480.393492, 87.486)
    7  = @(0.01287015847688175, 0.9961701388888889, 0.04085544673555932, 0.9893611111111111, 482.989044, 86.685)
    8  = @(0.01280095153107175, 0.9961994047619047, 0.04014146601044533, 0.9896785714285715, 487.935808, 87.48)
    9  = @(0.01280084208583996, 0.9961979166666666, 0.03979258810414816, 0.9898125, 482.874725, 86.952)
    10 = @(0.01272892271442571, 0.9962291666666666, 0.03870231707399497, 0.9899907407407407, 490.74664, 85.93899999999999)
    11 = @(0.012552012262605, 0.9962583333333332, 0.03807468705323368, 0.9901416666666666, 491.263767, 86.85899999999999)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
